$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.017617513053672
$ws.Range("D2").Value = 1.022783080465256
$ws.Range("E2").Value = 1.018994610517231
$ws.Range("I2").Value = 1.026783235987421
$ws.Range("J2").Value = 1.022830608662421
$ws.Range("K2").Value = 1.025616445634108
$ws.Range("L2").Value = 1.021839173692664
$ws.Range("N2").Value = 1.011711321753229
$ws.Range("B3").Value = 1.019999999999999
$ws.Range("C3").Value = 1.018673465791158
$ws.Range("D3").Value = 1.023536255111992
$ws.Range("E3").Value = 1.019891707136632
$ws.Range("I3").Value = 1.026961500359901
$ws.Range("J3").Value = 1.023521967413908
$ws.Range("K3").Value = 1.026176506126924
$ws.Range("L3").Value = 1.022541949170667
$ws.Range("N3").Value = 1.011939700539833
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.01935628074051
$ws.Range("D4").Value = 1.024022746659388
$ws.Range("E4").Value = 1.020472214266565
$ws.Range("I4").Value = 1.027074665774232
$ws.Range("J4").Value = 1.023968341123997
$ws.Range("K4").Value = 1.026537382793471
$ws.Range("L4").Value = 1.022996091964365
$ws.Range("N4").Value = 1.012087124210511
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.01964322703432
$ws.Range("D5").Value = 1.024227060407146
$ws.Range("E5").Value = 1.020716265174386
$ws.Range("I5").Value = 1.0271217172863
$ws.Range("J5").Value = 1.02415576135278
$ws.Range("K5").Value = 1.026688730941795
$ws.Range("L5").Value = 1.02318686965911
$ws.Range("N5").Value = 1.012149016577574
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.019691400190348
$ws.Range("D6").Value = 1.024261353406518
$ws.Range("E6").Value = 1.020757242717979
$ws.Range("I6").Value = 1.027129586745684
$ws.Range("J6").Value = 1.024187216228819
$ws.Range("K6").Value = 1.026714121577462
$ws.Range("L6").Value = 1.023218893623843
$ws.Range("N6").Value = 1.012159403617284
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.019360115359291
$ws.Range("D7").Value = 1.024025477525194
$ws.Range("E7").Value = 1.02047547526486
$ws.Range("I7").Value = 1.027075296535397
$ws.Range("J7").Value = 1.023970846367125
$ws.Range("K7").Value = 1.026539406547389
$ws.Range("L7").Value = 1.022998641711242
$ws.Range("N7").Value = 1.0120879515517
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.017974472060203
$ws.Range("D8").Value = 1.023037797574638
$ws.Range("E8").Value = 1.019297782960806
$ws.Range("I8").Value = 1.026843933100528
$ws.Range("J8").Value = 1.023064459681765
$ws.Range("K8").Value = 1.025806035095846
$ws.Range("L8").Value = 1.022076803992945
$ws.Range("N8").Value = 1.011788576390626
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.015529275598634
$ws.Range("D9").Value = 1.021290802736372
$ws.Range("E9").Value = 1.017222752757756
$ws.Range("I9").Value = 1.026419538257652
$ws.Range("J9").Value = 1.021459783479887
$ws.Range("K9").Value = 1.024502122441392
$ws.Range("L9").Value = 1.020447826147518
$ws.Range("N9").Value = 1.011258342680251
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.013896755377923
$ws.Range("D10").Value = 1.020121757219253
$ws.Range("E10").Value = 1.015839564047093
$ws.Range("I10").Value = 1.026125404197906
$ws.Range("J10").Value = 1.020384963726528
$ws.Range("K10").Value = 1.023625062954073
$ws.Range("L10").Value = 1.019358774494864
$ws.Range("N10").Value = 1.010903045832705
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.013189282522085
$ws.Range("D11").Value = 1.019614515522703
$ws.Range("E11").Value = 1.015240670516548
$ws.Range("I11").Value = 1.025995388336931
$ws.Range("J11").Value = 1.019918361966147
$ws.Range("K11").Value = 1.023243444978351
$ws.Range("L11").Value = 1.018886477858665
$ws.Range("N11").Value = 1.010748770377534
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.012926407609301
$ws.Range("D12").Value = 1.0194259480143
$ws.Range("E12").Value = 1.015018220410156
$ws.Range("I12").Value = 1.025946696259802
$ws.Range("J12").Value = 1.019744865348932
$ws.Range("K12").Value = 1.023101418091994
$ws.Range("L12").Value = 1.018710936230656
$ws.Range("N12").Value = 1.010691401087532
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.012982799164208
$ws.Range("D13").Value = 1.019466403407645
$ws.Range("E13").Value = 1.015065936456463
$ws.Range("I13").Value = 1.025957158898974
$ws.Range("J13").Value = 1.0197820891055
$ws.Range("K13").Value = 1.023131895870814
$ws.Range("L13").Value = 1.018748595465838
$ws.Range("N13").Value = 1.010703709914191
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.013167554999106
$ws.Range("D14").Value = 1.019598931640151
$ws.Range("E14").Value = 1.015222282605888
$ws.Range("I14").Value = 1.025991371560036
$ws.Range("J14").Value = 1.0199040243526
$ws.Range("K14").Value = 1.023231710640351
$ws.Range("L14").Value = 1.018871969769017
$ws.Range("N14").Value = 1.010744029530136
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.013281377484603
$ws.Range("D15").Value = 1.019680566085951
$ws.Range("E15").Value = 1.01531861337752
$ws.Range("I15").Value = 1.026012398327462
$ws.Range("J15").Value = 1.019979128837512
$ws.Range("K15").Value = 1.023293173119839
$ws.Range("L15").Value = 1.018947970223828
$ws.Range("N15").Value = 1.010768863227721
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.01394369572183
$ws.Range("D16").Value = 1.020155399367367
$ws.Range("E16").Value = 1.015879311417341
$ws.Range("I16").Value = 1.026133977035641
$ws.Range("J16").Value = 1.020415905314029
$ws.Range("K16").Value = 1.023650350825109
$ws.Range("L16").Value = 1.019390103890344
$ws.Range("N16").Value = 1.010913275537874
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.014358994437534
$ws.Range("D17").Value = 1.020452972377234
$ws.Range("E17").Value = 1.016231032323043
$ws.Range("I17").Value = 1.026209529857278
$ws.Range("J17").Value = 1.020689562996606
$ws.Range("K17").Value = 1.023873905088821
$ws.Range("L17").Value = 1.019667247239037
$ws.Range("N17").Value = 1.011003746571987
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.014601175160371
$ws.Range("D18").Value = 1.020626441566153
$ws.Range("E18").Value = 1.016436188669783
$ws.Range("I18").Value = 1.026253342459021
$ws.Range("J18").Value = 1.020849067331857
$ws.Range("K18").Value = 1.024004122433621
$ws.Range("L18").Value = 1.019828829819604
$ws.Range("N18").Value = 1.011056475346464
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.01468374305325
$ws.Range("D19").Value = 1.020685573080916
$ws.Range("E19").Value = 1.016506142269703
$ws.Range("I19").Value = 1.026268237978448
$ws.Range("J19").Value = 1.020903434587023
$ws.Range("K19").Value = 1.024048492969857
$ws.Range("L19").Value = 1.01988391335313
$ws.Range("N19").Value = 1.011074447464716
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.014314442612227
$ws.Range("D20").Value = 1.020421055961026
$ws.Range("E20").Value = 1.016193295635888
$ws.Range("I20").Value = 1.026201450239809
$ws.Range("J20").Value = 1.020660214070037
$ws.Range("K20").Value = 1.023849938249009
$ws.Range("L20").Value = 1.019637519666557
$ws.Range("N20").Value = 1.010994044173297
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.013113151445633
$ws.Range("D21").Value = 1.019559909663692
$ws.Range("E21").Value = 1.01517624240542
$ws.Range("I21").Value = 1.025981307779303
$ws.Range("J21").Value = 1.019868122429026
$ws.Range("K21").Value = 1.023202325322354
$ws.Range("L21").Value = 1.018835642137614
$ws.Range("N21").Value = 1.010732158189002
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.01235734253497
$ws.Range("D22").Value = 1.019017574185081
$ws.Range("E22").Value = 1.014536813240581
$ws.Range("I22").Value = 1.025840590622171
$ws.Range("J22").Value = 1.019369061842089
$ws.Range("K22").Value = 1.022793542985635
$ws.Range("L22").Value = 1.018330835654564
$ws.Range("N22").Value = 1.010567126810109
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01275805967831
$ws.Range("D23").Value = 1.019305161506096
$ws.Range("E23").Value = 1.014875783601223
$ws.Range("I23").Value = 1.025915405820054
$ws.Range("J23").Value = 1.019633721994025
$ws.Range("K23").Value = 1.023010398035835
$ws.Range("L23").Value = 1.018598503225723
$ws.Range("N23").Value = 1.010654648441559
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.014334573837054
$ws.Range("D24").Value = 1.020435477924933
$ws.Range("E24").Value = 1.016210347209535
$ws.Range("I24").Value = 1.026205101862361
$ws.Range("J24").Value = 1.020673475943551
$ws.Range("K24").Value = 1.023860768382394
$ws.Range("L24").Value = 1.019650952494483
$ws.Range("N24").Value = 1.010998428397653
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.016161836714717
$ws.Range("D25").Value = 1.021743218010405
$ws.Range("E25").Value = 1.017759170241455
$ws.Range("I25").Value = 1.026531231165477
$ws.Range("J25").Value = 1.02187551911245
$ws.Range("K25").Value = 1.024840589103763
$ws.Range("L25").Value = 1.020869497445361
$ws.Range("N25").Value = 1.011395739790958
